$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "Profiling?" -> "Profiling", highlighted yellow
$ws.Range("E5").Value = "Profiling"
$ws.Range("E5").Interior.Color = 65535

# Row 6: remove the "Profiling?" note entirely
$ws.Range("E6").ClearContents()

# Row 19: clear the "Profiling" note and remove its yellow highlight
$ws.Range("E19").ClearContents()
$ws.Range("E19").Interior.ColorIndex = -4142

# Row 22 & 23: remove the "Profiling" note
$ws.Range("E22").ClearContents()
$ws.Range("E23").ClearContents()

# Row 25: clear the "Profiling" note and remove its yellow highlight
$ws.Range("E25").ClearContents()
$ws.Range("E25").Interior.ColorIndex = -4142

# Row 29: keep the "Profiling" note but highlight it yellow
$ws.Range("E29").Interior.Color = 65535

# Row 30: remove the "Profiling" note entirely
$ws.Range("E30").ClearContents()

$ws.Range("G11").Select()
